# Multiply "actualEmissionChangePercent" (column V) by 100 for every data row,
# so the backend now stores the value already scaled to percent instead of a
# fraction (row 1 is the header "actualEmissionChangePercent" and is left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$colV = 22  # column V

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colV)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current * 100
    }
}
